# Items.xlsx — "add edit and delete feature"
# Appends 6 new rows (8-13) of item data to the existing "Items" sheet,
# matching the shape of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value into a cell as literal TEXT (shared string), even
# when it looks like a number (e.g. "12", "56", "44565655"), then strip the
# quote-prefix/number-format bookkeeping that Excel adds for apostrophe
# entry so the cell is left with no explicit style — matching the plain
# t="s" cells already used for the "Item"/"Finsih"/"Rate-SQFT" columns in
# this sheet.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# Helper: set a plain numeric cell.
function Set-NumberValue {
    param($range, $number)
    $range.Value = $number
}

# Row data: A, B(text), C, D, E, F, G, H(text), I(text), J, K
$rows = @(
    @{ Row=8;  A=7;  B="car case";    C=18;     D=664;    E=540;     F=2;      G=7.719;             H="WITH GROOVE";                  I="12";     J=280;    K=2161.32 },
    @{ Row=9;  A=8;  B="Www";         C=12;     D=550;    E=700;     F=3;      G=12.432;            H="Sss";                          I="56";     J=256;    K=3182.592 },
    @{ Row=10; A=9;  B="Www";         C=12;     D=550;    E=700;     F=3;      G=12.432;            H="Sss";                          I="56";     J=256;    K=3182.592 },
    @{ Row=11; A=10; B="Sss";         C=12;     D=550;    E=600;     F=3;      G=10.656;            H="Ass";                          I="12";     J=300;    K=3196.8 },
    @{ Row=12; A=11; B="44565655";    C=123456; D=123456; E=123456;  F=123456; G=20253807542.593;   H="ggrgr";                        I="123456"; J=123456; K=2500454063978361 },
    @{ Row=13; A=12; B="new carcase"; C=444;    D=444;    E=4435345; F=435343; G=9228094145.572;    H="WITH GROOVE plus with gropve"; I="acrylic"; J=3333;  K=30757237787191.477 }
)

foreach ($r in $rows) {
    $n = $r.Row

    Set-NumberValue $ws.Range("A$n") $r.A
    Set-TextValue   $ws.Range("B$n") $r.B
    Set-NumberValue $ws.Range("C$n") $r.C
    Set-NumberValue $ws.Range("D$n") $r.D
    Set-NumberValue $ws.Range("E$n") $r.E
    Set-NumberValue $ws.Range("F$n") $r.F
    Set-NumberValue $ws.Range("G$n") $r.G
    Set-TextValue   $ws.Range("H$n") $r.H
    Set-TextValue   $ws.Range("I$n") $r.I
    Set-NumberValue $ws.Range("J$n") $r.J
    Set-NumberValue $ws.Range("K$n") $r.K
}
